$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("How many curve shades can I create?", "llama3.2:latest", "According to the document, you can create 250 curve shades."),
    @("What's the maximum number of lithology types in an log?", "llama3.2:latest", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("How many tracks can you define in one ODF?", "llama3.2:latest", "According to the Track Settings dialog box, the number of tracks that can be defined is 200."),
    @("How many tracks can you define in one ODF?", "llama3.2:latest", "According to the Track Settings dialog box, the number of tracks that can be defined is 200."),
    @("How many curve shades can I create?", "llama3.2:latest", "According to the document, you can create 250 curve shades."),
    @("What's the maximum number of lithology types in an log?", "llama3.2:latest", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("How many tracks can you define in one ODF?", "llama3.2:latest", "According to the Track Settings dialog box, the number of tracks that can be defined is 200."),
    @("How many curve shades can I create?", "llama3.2:latest", "According to the document, you can create 250 curve shades per plot."),
    @("What's the maximum number of lithology types in an log?", "llama3.2:latest", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("How many tracks can you define in one ODF?", "llama3.2:latest", "According to the Track Settings dialog box, the number of tracks that can be defined is 200."),
    @("How many curve shades can I create?", "llama3.2:latest", "According to the document, you can create 250 curve shades per plot."),
    @("What's the maximum number of lithology types in an log?", "llama3.2:latest", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("How many tracks can you define in one ODF?", "llama3.2:latest", "According to the Track Settings dialog box, the number of tracks that can be defined is 200."),
    @("How many curve shades can I create?", "llama3.2:latest", "According to the document, you can create 250 curve shades per plot."),
    @("What's the maximum number of lithology types in an log?", "llama3.2:latest", "The maximum number of lithology types that can be recorded in a log is 450."),
    @("How many tracks can you define in one ODF?", "llama3.2:latest", "According to the Track Settings dialog box, the number of tracks that can be defined is 200."),
    @("How many curve shades can I create?", "llama3.2:latest", "According to the document, you can create 250 curve shades per plot.")
)

$startRow = 43
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

$wb.Save()
